$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Sincerely," (the sign-off line,
# right after the "Let's check now" paragraph and right before the
# signature-name paragraph).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Sincerely,")) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    Write-Host "Could not find the 'Sincerely,' paragraph"
} else {
    $target = $d.Paragraphs.Item($targetIndex)

    # Insert the two replacement paragraphs right after the current one
    # *before* we touch its text/formatting, so the new paragraphs don't
    # inherit the magenta highlight we are about to apply.
    $target.Range.InsertParagraphAfter()

    $blank = $d.Paragraphs.Item($targetIndex + 1)
    $blank.Range.InsertParagraphAfter()

    $closing = $d.Paragraphs.Item($targetIndex + 2)
    $closing.Range.Text = "Sincerely,"

    # Now turn the original paragraph into the new "Verifying the change
    # update" line, highlighted in magenta.
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.Text = "Verifying the change update"
    $target.Range.HighlightColorIndex = 5
}
